$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.505.72"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.623.80"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'211.69"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'23.25"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("D10").Value = "'0.0610"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'0.0880"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "1.852.67"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.621.10"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "'65.37"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "27.493.60"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'229.49"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'10.45"
$ws.Range("E22").Value = "  +4.12%  "
$ws.Range("D23").Value = "'4.35"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("E24").Value = "  +8.41%  "
$ws.Range("D25").Value = "'149.08"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "'6.88"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'15.50"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "'0.0483"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "1.468.83"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").Value = "'0.943"
$ws.Range("E37").Value = "  +4.25%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").Value = "'1.03"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'67.77"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").Value = "'2.48"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").Value = "'5.32"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("D47").Value = "'1.77"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").Value = "1.763.31"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").Value = "'87.21"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").Value = "'0.0993"
$ws.Range("E51").Value = "  +0.55%  "
